$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the average_doctor / average_doctor_old header labels
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Update the *_old metric columns (Ada_old, Avey_old, Buoy_old, K health_old,
# WebMD_old, doctor_MA_old, doctor_NJ_old, doctor_TH_old) and the average_doctor /
# average_doctor_old data columns with the newly computed "harvard case classification" values
$values = @{
    "AI4" = 0.195
    "AJ4" = 0.064
    "AK4" = 0.254
    "AU4" = 0.144
    "AV4" = 0.027
    "AW4" = 0.164
    "BA4" = 1.974
    "BB4" = 0.165
    "BC4" = 0.407
    "BG4" = 0.74
    "BH4" = 0.135
    "BI4" = 0.368
    "BM4" = 0.696
    "BN4" = 0.08500000000000001
    "BO4" = 0.291
    "BP4" = 0.658
    "BQ4" = 0.66
    "E4" = 0.383
    "F4" = 0.08500000000000001
    "G4" = 0.292
    "N4" = 0.41
    "O4" = 0.064
    "P4" = 0.252
    "W4" = 0.226
    "X4" = 0.104
    "Y4" = 0.322
    "AI5" = 0.238
    "AJ5" = 0.096
    "AK5" = 0.309
    "AU5" = 0.303
    "AV5" = 0.102
    "AW5" = 0.319
    "BA5" = 1.396
    "BB5" = 0.08599999999999999
    "BC5" = 0.293
    "BG5" = 0.42
    "BH5" = 0.051
    "BI5" = 0.225
    "BM5" = 0.584
    "BN5" = 0.074
    "BO5" = 0.272
    "BP5" = 0.465
    "BQ5" = 0.462
    "E5" = 0.503
    "F5" = 0.103
    "G5" = 0.32
    "N5" = 0.773
    "O5" = 0.075
    "P5" = 0.274
    "W5" = 0.234
    "X5" = 0.116
    "Y5" = 0.34
    "AI6" = 0.214
    "AU6" = 0.195
    "BA6" = 1.625
    "BG6" = 0.536
    "BM6" = 0.635
    "BP6" = 0.542
    "BQ6" = 0.54
    "E6" = 0.435
    "N6" = 0.536
    "W6" = 0.23
    "AI7" = 0.228
    "AU7" = 0.248
    "BA7" = 1.477
    "BG7" = 0.46
    "BM7" = 0.603
    "BP7" = 0.492
    "BQ7" = 0.49
    "E7" = 0.473
    "N7" = 0.657
    "W7" = 0.232
    "AI8" = 0.22
    "AJ8" = 0.095
    "AK8" = 0.308
    "AU8" = 0.245
    "AV8" = 0.077
    "AW8" = 0.278
    "BA8" = 1.758
    "BB8" = 0.136
    "BC8" = 0.369
    "BG8" = 0.582
    "BH8" = 0.101
    "BI8" = 0.318
    "BM8" = 0.715
    "BN8" = 0.064
    "BO8" = 0.254
    "BP8" = 0.586
    "BQ8" = 0.59
    "E8" = 0.52
    "F8" = 0.131
    "G8" = 0.361
    "N8" = 0.772
    "O8" = 0.056
    "P8" = 0.237
    "W8" = 0.227
    "X8" = 0.108
    "Y8" = 0.328
    "AI9" = 0.125
    "AJ9" = 0.109
    "AK9" = 0.331
    "BA9" = 1.73
    "BB9" = 0.246
    "BC9" = 0.496
    "BG9" = 0.625
    "BH9" = 0.234
    "BI9" = 0.484
    "BM9" = 0.667
    "BN9" = 0.222
    "BO9" = 0.471
    "BP9" = 0.577
    "BQ9" = 0.5629999999999999
    "E9" = 0.438
    "F9" = 0.246
    "G9" = 0.496
    "N9" = 0.646
    "O9" = 0.229
    "P9" = 0.478
    "W9" = 0.125
    "X9" = 0.109
    "Y9" = 0.331
    "AI10" = 0.25
    "AJ10" = 0.188
    "AK10" = 0.433
    "AU10" = 0.229
    "AV10" = 0.177
    "AW10" = 0.42
    "BA10" = 2.042
    "BG10" = 0.667
    "BH10" = 0.222
    "BI10" = 0.471
    "BM10" = 0.875
    "BN10" = 0.109
    "BO10" = 0.331
    "BP10" = 0.681
    "BQ10" = 0.695
    "E10" = 0.5620000000000001
    "F10" = 0.246
    "G10" = 0.496
    "N10" = 0.854
    "O10" = 0.125
    "P10" = 0.353
    "W10" = 0.271
    "X10" = 0.197
    "Y10" = 0.444
    "AI11" = 0.25
    "AJ11" = 0.188
    "AK11" = 0.433
    "AU11" = 0.354
    "AV11" = 0.229
    "AW11" = 0.478
    "BA11" = 2.042
    "BG11" = 0.667
    "BH11" = 0.222
    "BI11" = 0.471
    "BM11" = 0.875
    "BN11" = 0.109
    "BO11" = 0.331
    "BP11" = 0.681
    "BQ11" = 0.695
    "E11" = 0.583
    "F11" = 0.243
    "G11" = 0.493
    "N11" = 0.896
    "O11" = 0.093
    "P11" = 0.305
    "W11" = 0.271
    "X11" = 0.197
    "Y11" = 0.444
    "AI12" = 1.917
    "AJ12" = 0.91
    "AK12" = 0.954
    "AU12" = 2.647
    "AV12" = 1.758
    "AW12" = 1.326
    "BA12" = 3.603
    "BB12" = 0.332
    "BC12" = 0.576
    "BG12" = 1.062
    "BH12" = 0.059
    "BI12" = 0.242
    "BM12" = 1.333
    "BN12" = 0.413
    "BO12" = 0.642
    "BP12" = 1.201
    "BQ12" = 1.257
    "E12" = 1.429
    "F12" = 0.673
    "G12" = 0.821
    "N12" = 1.667
    "O12" = 1.644
    "P12" = 1.282
    "W12" = 1.846
    "X12" = 0.746
    "Y12" = 0.863
    "AI13" = 1.4
    "AJ13" = 0.398
    "AK13" = 0.631
    "AU13" = 2.523
    "AV13" = 1.319
    "AW13" = 1.148
    "BA13" = 2.541
    "BB13" = 0.319
    "BC13" = 0.5649999999999999
    "BG13" = 0.633
    "BH13" = 0.08699999999999999
    "BI13" = 0.294
    "BM13" = 0.999
    "BN13" = 0.378
    "BO13" = 0.615
    "BP13" = 0.847
    "BQ13" = 0.794
    "E13" = 1.765
    "F13" = 0.922
    "G13" = 0.96
    "N13" = 2.352
    "O13" = 1.21
    "P13" = 1.1
    "W13" = 1.11
    "X13" = 0.175
    "Y13" = 0.419
}

foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}
